$d = $word.ActiveDocument

# Locate the paragraph that contains the ellipsis "…" -- the new paragraph
# needs to be added right after it, and right before the trailing empty
# paragraph that precedes the section break.
$ellipsis = [char]0x2026
$ellipsisPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains($ellipsis)) {
        $ellipsisPara = $p
    }
}

# Insert a brand new (empty) paragraph right after it, then fill in its text.
$ellipsisPara.Range.InsertParagraphAfter()
$newRange = $ellipsisPara.Next().Range
$newRange.Text = "Adding some details about the version management, that is helping us to manage and control our code in the best way possible, we can do it using various tools today we are exploring it using github."
